$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A-D (Sending cluster / Ligand symbol / Receptor symbol / Target cluster)
# already hold the correct text for rows 2-4 (ECs / Il1b / Il1r1 / {ECs,FAPs,MuSCs});
# only the numeric NATMI metrics (re-run with updated TPM data) and the row count change.

# Row 2 (ECs -> Il1b -> Il1r1 -> ECs)
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1046376666666667
$ws.Range("H2").Value = 0.313913
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("M2").Value = 5.828378333333333
$ws.Range("N2").Value = 17.485135
$ws.Range("O2").Value = 0.1413867973615592
$ws.Range("P2").Value = 0.1413867973615592
$ws.Range("Q2").Value = 0.6098679092505556
$ws.Range("R2").Value = 5.488811183255
$ws.Range("S2").Value = 0.1413867973615592
$ws.Range("T2").Value = 0.1413867973615592

# Row 3 (ECs -> Il1b -> Il1r1 -> FAPs)
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1046376666666667
$ws.Range("H3").Value = 0.313913
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("O3").Value = 0.604557320991465
$ws.Range("P3").Value = 0.604557320991465
$ws.Range("Q3").Value = 2.607740724420889
$ws.Range("R3").Value = 23.469666519788
$ws.Range("S3").Value = 0.604557320991465
$ws.Range("T3").Value = 0.604557320991465

# Row 4 (ECs -> Il1b -> Il1r1 -> MuSCs)
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1046376666666667
$ws.Range("H4").Value = 0.313913
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("M4").Value = 10.47292833333333
$ws.Range("N4").Value = 31.418785
$ws.Range("O4").Value = 0.2540558816469758
$ws.Range("P4").Value = 0.2540558816469758
$ws.Range("Q4").Value = 1.095862783967222
$ws.Range("R4").Value = 9.862765055704999
$ws.Range("S4").Value = 0.2540558816469758
$ws.Range("T4").Value = 0.2540558816469758

# The old rows 5-7 (MuSCs as sender) are gone from the re-run output entirely.
$ws.Range("A5:T7").Delete() | Out-Null
